# Apply the "terminacion analisis simulacion base" edit:
# - Swap the model labels for rows 3 and 4 (DeepAR <-> AREPD)
# - Update the Sensibilidad_Media (B) and Std_entre_d (C) values for rows 2-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap labels in A3 and A4 (row 3 was "DeepAR", row 4 was "AREPD")
$ws.Range("A3").Value = "AREPD"
$ws.Range("A4").Value = "DeepAR"

# Update numeric values for B2:C10
$ws.Range("B2").Value = 2999369305732183
$ws.Range("C2").Value = 6297321029801379

$ws.Range("B3").Value = 2957864556379087
$ws.Range("C3").Value = 6210180144467948

$ws.Range("B4").Value = 2947246238922365
$ws.Range("C4").Value = 6187887764827355

$ws.Range("B5").Value = 2140874874689569
$ws.Range("C5").Value = 4494874881394919

$ws.Range("B6").Value = 1312607574687019
$ws.Range("C6").Value = 2755887130005798

$ws.Range("B7").Value = 1227315061997159
$ws.Range("C7").Value = 2576811884227306

$ws.Range("B8").Value = 922533738916231.6
$ws.Range("C8").Value = 1936907386082162

$ws.Range("B9").Value = 150416090478961.8
$ws.Range("C9").Value = 315806398270839.8

$ws.Range("B10").Value = 164569420.8808171
$ws.Range("C10").Value = 345521813.0988252
